# Import guru dan karawan: add a "NOMER TELEPON" (phone number) column to
# the teacher-import template, between NIK (G) and AGAMA (old H, now I).
#
# Inserting a whole column shifts AGAMA/ALAMAT (and everything to their
# right) one column to the right, which also carries along their data
# validations (religion dropdown) and column widths automatically - just
# like using Excel's "Insert Sheet Columns" on column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H, pushing AGAMA/ALAMAT/etc. one column right.
$ws.Columns("H:H").Insert() | Out-Null

# New header (row 1) and example value (row 2) for the phone number column.
$ws.Range("H1").Value = "NOMER TELEPON"
$ws.Range("H2").Value = "0878 9878 7878"

# Match the author's final cursor position.
$ws.Range("E24").Select() | Out-Null
